$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7947
$ws1.Range("F5").Value = 5803
$ws1.Range("F6").Value = 488
$ws1.Range("F7").Value = 83
$ws1.Range("F10").Value = 276
$ws1.Range("F11").Value = 342

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 87

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7947
$ws4.Range("F5").Value = 5803
$ws4.Range("F6").Value = 488
$ws4.Range("F7").Value = 83
$ws4.Range("F10").Value = 276
$ws4.Range("F11").Value = 87
$ws4.Range("F14").Value = 342
